$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '29.038.10'
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = '  -0.45%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.830.36'
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = '  -0.27%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9984'
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = '  -0.08%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '241.24'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  -0.35%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.6271'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  -4.77%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.9998'
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  -0.03%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.07585'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  +2.19%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '44.92'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  +7.38%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.2914'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  -0.55%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '22.79'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  -0.63%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.07644'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  -1.66%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.829.85'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  -2.32%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '4.957'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  -0.48%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.6653'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  +0.04%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '82.31'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  -0.60%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.000009115'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  +6.14%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '5.986'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  -2.01%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '29.019.50'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  -0.56%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '224.95'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  -0.78%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '12.33'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  -0.93%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.9999'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  -0.07%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.200'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  +1.46%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.9999'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  +0.00%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '159.90'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  +0.31%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '8.420'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  -2.08%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.1362'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  -2.27%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '17.83'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  -0.55%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.498'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  -1.27%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.033'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  -0.25%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.047'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  -1.55%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.203'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  +0.87%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.05196'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  -1.42%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.845'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  -1.09%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.153'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  +0.69%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.7322'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  -0.91%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.613'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  -1.67%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.281.34'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  -1.53%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.761'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  +0.92%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '6.468'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  +7.18%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.8906'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  -3.17%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.000'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  +0.07%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '101.67'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  -0.97%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.977.91'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  -2.76%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.5108'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  -0.68%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '63.84'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  +0.51%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.3980'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  -0.54%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.07214'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  -16.02%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '8.822'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  +0.70%  '
